# Update cryptocurrency price / volume data as scraped by GitHub Actions.
# Column D holds text-like price strings (dotted thousands separators), so
# force Text number format before assigning to avoid Excel auto-converting
# them into numeric values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.838.01"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.06"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.32"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5077"
$ws.Range("E7").Value = "  -2.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3651"
$ws.Range("E8").Value = "  -2.34%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07156"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8883"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.51"
$ws.Range("E11").Value = "  -1.26%  "

# Rows 12 and 13 swap rank (TRON moves above WrappedEther) with updated values.
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07481"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.873.58"
$ws.Range("E13").Value = "  +0.11%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.06"
$ws.Range("E14").Value = "  +3.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.210"
$ws.Range("E15").Value = "  -1.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008479"
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.13"
$ws.Range("E18").Value = "  +0.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.861.68"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.987"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.123.11"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.32"
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.348"
$ws.Range("E24").Value = "  -1.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.95"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.769"
$ws.Range("E26").Value = "  -3.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.84"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.089"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.47"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.670"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.699"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09115"
$ws.Range("E32").Value = "  -1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05026"
$ws.Range("E33").Value = "  -2.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7467"
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.948"
$ws.Range("E35").Value = "  -4.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.148"
$ws.Range("E36").Value = "  -1.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.205"
$ws.Range("E37").Value = "  +2.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.495"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01976"
$ws.Range("E39").Value = "  -2.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5532"
$ws.Range("E40").Value = "  +4.02%  "
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.564"
$ws.Range("E42").Value = "  +0.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "115.85"
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.558"
$ws.Range("E44").Value = "  +2.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1481"
$ws.Range("E45").Value = "  +0.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4731"
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.00"
$ws.Range("E48").Value = "  -0.03%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.03"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.549"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.78"
$ws.Range("E51").Value = "  -1.56%  "
